# The "missing_data" worksheet dropped two samples (ID "RM 232" and ID "SC 92")
# from the data table. Removing these two rows shifts all subsequent rows up by
# two, shrinking the used range from A1:F35 down to A1:F33. Because column B
# ("A") represents randomly-missing measurements, the set of rows left blank in
# column B also changed for the rows that shifted into new positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "RM 232" row (original row 26).
$ws.Rows.Item(26).Delete()

# Remove the "SC 92" row (originally row 28, now row 27 after the delete above).
$ws.Rows.Item(27).Delete()

# After the deletions the rows have shifted as follows (old -> new):
#   27 (SC 5)   -> 26
#   29 (SC 101) -> 27
#   30 (SC 105) -> 28
#   31 (SC 119) -> 29
#   32 (SC 120) -> 30
#   33 (SC 132) -> 31
#   34 (SC 193) -> 32
#   35 (SC 232) -> 33
#
# Update column B ("A") so the missing-value pattern matches the new row order:
# "SC 5" (row 26) and "SC 119" (row 29) are now the missing entries, while
# "SC 101" (row 27) now has its value filled in.
$ws.Range("B26").ClearContents()
$ws.Range("B27").Value = -20.4
$ws.Range("B29").ClearContents()
